$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in "reference product" values for each of the four activity blocks,
# both in the activity header rows (column B) and in the production
# exchange rows (column B of the exchange table).

# Egg block (header row 4, exchange row 11)
$ws.Range("B4").Value = "egg"
$ws.Range("B11").Value = "egg"

# Sea bass or sea bream block (header row 15, exchange row 22)
$ws.Range("B15").Value = "sea bass or sea bream"
$ws.Range("B22").Value = "sea bass or sea bream"
$ws.Range("B22").Style = "Normal"

# Large trout block (header row 26, exchange row 33)
$ws.Range("B26").Value = "large trout"
$ws.Range("B33").Value = "large trout"

# Small trout block (header row 37, exchange row 44)
$ws.Range("B37").Value = "small trout"
$ws.Range("B44").Value = "small trout"

# Update the selected cell, matching the saved cursor position in the file.
$ws.Range("B49").Select()
